$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the ID and Name for the first employee row
$ws.Range("C2:C7").Value = "EN-4-072"
$ws.Range("D2").Value = "Agus C"

# Update the active selection to D9 (matches sheetView selection in diff)
$ws.Range("D9").Select()
